$d = $word.ActiveDocument

$section = $d.Sections.Item(1)
$header = $section.Headers.Item(1)  # wdHeaderFooterPrimary = 1

$range = $header.Range
$range.InsertAfter("Questionnaire 31")
$range.set_Style("Header")
$range.ParagraphFormat.Alignment = 1  # wdAlignParagraphCenter

# Apply the run-level font formatting to the text only (not the paragraph
# mark) so the w:rPr stays on the <w:r>, matching how Word records a
# formatted run rather than stamping the paragraph mark too.
$textRange = $range.Duplicate
[void]$textRange.MoveEnd(1, -1)  # wdCharacter = 1
$textRange.Font.Name = "Arial"
$textRange.Font.Size = 12

$lastChar = $range.Duplicate
$lastChar.SetRange($range.End - 1, $range.End)
$lastChar.Font.Name = "Arial"
$lastChar.Font.Size = 12
